$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.982.03'
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("E3").Value = '  -0.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.956.91'
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = '  -0.19%  '

# Row 5
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.09'
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("E6").Value = '  -0.02%  '

# Row 7
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4852'
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2943'
$ws.Range("D8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = '  +4.67%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07120'
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("E10").Value = '  +2.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.62'
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = '  +1.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '107.55'
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.961.29'
$ws.Range("D12").Style = "Normal"

# Row 13
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07774'
$ws.Range("D13").Style = "Normal"

# Row 14
$ws.Range("E14").Value = '  -0.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.380'
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = '  +0.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7050'
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("E16").Value = '  -2.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '278.60'
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.000.36'
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("E18").Value = '  +1.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007846'
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = '  +1.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.34'
$ws.Range("D19").Style = "Normal"

# Row 20
$ws.Range("E20").Value = '  -0.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.213.71'
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = '  +0.01%  '

# Row 22
$ws.Range("E22").Value = '  -1.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.510'
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = '  -0.31%  '

# Row 24
$ws.Range("E24").Value = '  -1.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.516'
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = '  -2.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.762'
$ws.Range("D25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.47'
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = '  -1.51%  '

# Row 28
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.175'
$ws.Range("D28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = '  -1.23%  '

# Row 30
$ws.Range("E30").Value = '  -3.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.402'
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("E31").Value = '  -1.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.570'
$ws.Range("D31").Style = "Normal"

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("E32").Value = '  -4.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.604'
$ws.Range("D32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.421'
$ws.Range("D33").Style = "Normal"

# Row 34
$ws.Range("E34").Value = '  -3.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04912'
$ws.Range("D34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = '  -2.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7543'
$ws.Range("D35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.169'
$ws.Range("D36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.731'
$ws.Range("D37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = '  -1.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02007'
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("E39").Value = '  -1.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.678'
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.513'
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E41").Value = '  +8.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '77.83'
$ws.Range("D41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.118'
$ws.Range("D42").Style = "Normal"

# Row 43
$ws.Range("E43").Value = '  +0.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8963'
$ws.Range("D43").Style = "Normal"

# Row 44
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '109.58'
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4455'
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = '  +5.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.874'
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = '  +2.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '981.91'
$ws.Range("D48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = '  -1.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1251'
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = '  -1.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.280'
$ws.Range("D50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = '  +0.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.00'
$ws.Range("D51").Style = "Normal"
